# Update RAF values for hydrogen plants (RAF-capacity sheet) and add an
# explanatory note about the change on the About sheet.

$wb = $excel.ActiveWorkbook

# --- RAF-capacity: lower the capacity-credit RAF for hydrogen plants ---
# Row 24 = "hydrogen combustion turbine", Row 25 = "hydrogen combined cycle"
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 0.1
$wsCap.Range("B25").Value = 0.1

# Touch/activate this sheet so the selection state is recorded, then move on
# (About will be (re)activated afterwards so it ends up as the visible tab).
$wsCap.Activate()
$wsCap.Range("B26").Select()

# --- About: document why the hydrogen RAF credit was lowered ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A45").Value = "We also apply a very low credit for hydrogen plants because they require "
$wsAbout.Range("A46").Value = "access to a hydrogen supplier, most likely pipeline deliery, which "
$wsAbout.Range("A47").Value = "does not exist in the US today and they would only be used in certain "
$wsAbout.Range("A48").Value = "unique circumstances."

$wsAbout.Activate()
$wsAbout.Range("A49").Select()
